$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 8333.5
$ws.Range("J18").Value = 6618
$ws.Range("L18").Value = 6618
$ws.Range("N18").Value = -7186
$ws.Range("H62").Value = 4786.6665
$ws.Range("I62").Value = 4767.273
$ws.Range("K62").Value = 4767.273
$ws.Range("M62").Value = -4143.273
$ws.Range("H65").Value = 4786.6665
$ws.Range("I65").Value = 4767.273
$ws.Range("K65").Value = 23836.365
$ws.Range("M65").Value = -20716.365
$ws.Range("H74").Value = 2859.8
$ws.Range("I74").Value = 2859.8
$ws.Range("K74").Value = 2859.8
$ws.Range("M74").Value = -1923.8
$ws.Range("H77").Value = 2859.8
$ws.Range("I77").Value = 2859.8
$ws.Range("K77").Value = 14299
$ws.Range("M77").Value = -9619
$ws.Range("H112").Value = 3488.5173
$ws.Range("J112").Value = 3985.4583
$ws.Range("L112").Value = 11956.3749
$ws.Range("N112").Value = -14172.3749
$ws.Range("H132").Value = 5045.718
$ws.Range("I132").Value = 5128.091
$ws.Range("K132").Value = 15384.273
$ws.Range("M132").Value = -12854.273
$ws.Range("H137").Value = 12405.08
$ws.Range("I137").Value = 15943.471
$ws.Range("J137").Value = 4886
$ws.Range("K137").Value = 47830.413
$ws.Range("L137").Value = 14658
$ws.Range("M137").Value = -45280.413
$ws.Range("N137").Value = -19758

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6007.1353
$ws.Range("I32").Value = 6007.1353
$ws.Range("K32").Value = 6007.1353
$ws.Range("M32").Value = -5720.1353
$ws.Range("H45").Value = 111509.1
$ws.Range("I45").Value = 148081.78
$ws.Range("K45").Value = 148081.78
$ws.Range("M45").Value = -147704.78
$ws.Range("H61").Value = 4322.1025
$ws.Range("I61").Value = 4317.303
$ws.Range("K61").Value = 4317.303
$ws.Range("M61").Value = -4105.303
$ws.Range("H74").Value = 3032.4
$ws.Range("I74").Value = 1340.8462
$ws.Range("K74").Value = 1340.8462
$ws.Range("M74").Value = -466.8462
$ws.Range("H77").Value = 3032.4
$ws.Range("I77").Value = 1340.8462
$ws.Range("K77").Value = 6704.231
$ws.Range("M77").Value = -2336.231
$ws.Range("H132").Value = 4088.1384
$ws.Range("I132").Value = 3795.6667
$ws.Range("K132").Value = 11387.0001
$ws.Range("M132").Value = -8857.000100000001
$ws.Range("H136").Value = 4322.1025
$ws.Range("I136").Value = 4317.303
$ws.Range("K136").Value = 12951.909
$ws.Range("M136").Value = -10401.909

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 736.7273
$ws.Range("I22").Value = 782.7143
$ws.Range("K22").Value = 782.7143
$ws.Range("M22").Value = -609.7143
$ws.Range("H134").Value = 6265.64
$ws.Range("I134").Value = 6758.909
$ws.Range("K134").Value = 20276.727
$ws.Range("M134").Value = -17741.727
$ws.Range("H138").Value = 94709.42999999999
$ws.Range("J138").Value = 94709.42999999999
$ws.Range("L138").Value = 94709.42999999999
$ws.Range("N138").Value = -104989.43

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2812.7778
$ws.Range("I134").Value = 2842.2
$ws.Range("K134").Value = 8526.599999999999
$ws.Range("M134").Value = -5991.599999999999
$ws.Range("H141").Value = 220780.83
$ws.Range("J141").Value = 306881.5
$ws.Range("L141").Value = 306881.5
$ws.Range("N141").Value = -317241.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 12055.454
$ws.Range("J81").Value = 13878.889
$ws.Range("L81").Value = 41636.667
$ws.Range("N81").Value = -43882.667
$ws.Range("H84").Value = 12055.454
$ws.Range("J84").Value = 13878.889
$ws.Range("L84").Value = 124910.001
$ws.Range("N84").Value = -136142.001
$ws.Range("H103").Value = 1249.1177
$ws.Range("I103").Value = 618.375
$ws.Range("J103").Value = 1809.7778
$ws.Range("K103").Value = 1855.125
$ws.Range("L103").Value = 5429.3334
$ws.Range("M103").Value = -976.125
$ws.Range("N103").Value = -7187.3334

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6523.816
$ws.Range("I122").Value = 3876.0908
$ws.Range("K122").Value = 11628.2724
$ws.Range("M122").Value = -9178.2724
$ws.Range("H132").Value = 1601.3334
$ws.Range("I132").Value = 1457.375
$ws.Range("K132").Value = 4372.125
$ws.Range("M132").Value = -1842.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13885.9375
$ws.Range("I22").Value = 29141.715
$ws.Range("J22").Value = 2020.3334
$ws.Range("K22").Value = 29141.715
$ws.Range("L22").Value = 2020.3334
$ws.Range("M22").Value = -28846.715
$ws.Range("N22").Value = -2610.3334
$ws.Range("H27").Value = 13885.9375
$ws.Range("I27").Value = 29141.715
$ws.Range("J27").Value = 2020.3334
$ws.Range("K27").Value = 29141.715
$ws.Range("L27").Value = 2020.3334
$ws.Range("M27").Value = -29034.715
$ws.Range("N27").Value = -2234.3334
$ws.Range("H68").Value = 3874.6875
$ws.Range("I68").Value = 2097
$ws.Range("J68").Value = 4467.25
$ws.Range("K68").Value = 2097
$ws.Range("L68").Value = 4467.25
$ws.Range("M68").Value = -1348
$ws.Range("N68").Value = -5965.25
$ws.Range("H71").Value = 3874.6875
$ws.Range("I71").Value = 2097
$ws.Range("J71").Value = 4467.25
$ws.Range("K71").Value = 10485
$ws.Range("L71").Value = 22336.25
$ws.Range("M71").Value = -6741
$ws.Range("N71").Value = -29824.25
$ws.Range("H100").Value = 4645.385
$ws.Range("I100").Value = 2880
$ws.Range("K100").Value = 2880
$ws.Range("M100").Value = -2339
$ws.Range("H122").Value = 6398
$ws.Range("I122").Value = 7824.2856
$ws.Range("K122").Value = 23472.8568
$ws.Range("M122").Value = -21022.8568
$ws.Range("H132").Value = 681871.9399999999
$ws.Range("I132").Value = 1149707.4
$ws.Range("J132").Value = 6109.6665
$ws.Range("K132").Value = 3449122.2
$ws.Range("L132").Value = 18328.9995
$ws.Range("M132").Value = -3446592.2
$ws.Range("N132").Value = -23388.9995
$ws.Range("H136").Value = 4640.425
$ws.Range("I136").Value = 2449
$ws.Range("J136").Value = 8710.214
$ws.Range("K136").Value = 7347
$ws.Range("L136").Value = 26130.642
$ws.Range("M136").Value = -4797
$ws.Range("N136").Value = -31230.642

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1539.6428
$ws.Range("I113").Value = 763.6842
$ws.Range("J113").Value = 3177.7778
$ws.Range("K113").Value = 2291.0526
$ws.Range("L113").Value = 9533.3334
$ws.Range("M113").Value = -121.0526
$ws.Range("N113").Value = -13873.3334
$ws.Range("H126").Value = 16578.586
$ws.Range("I126").Value = 22888.63
$ws.Range("J126").Value = 4589.5
$ws.Range("K126").Value = 68665.89
$ws.Range("L126").Value = 13768.5
$ws.Range("M126").Value = -66195.89
$ws.Range("N126").Value = -18708.5
$ws.Range("H132").Value = 9184.393
$ws.Range("I132").Value = 11068.324
$ws.Range("J132").Value = 4205.4287
$ws.Range("K132").Value = 33204.972
$ws.Range("L132").Value = 12616.2861
$ws.Range("M132").Value = -30674.972
$ws.Range("N132").Value = -17676.2861
